$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
Write-Host "before A16:" ($ws.Range("A16").Value2)
Write-Host "before A17:" ($ws.Range("A17").Value2)
$ws.Rows.Item(16).Delete()
Write-Host "after A16:" ($ws.Range("A16").Value2)
Write-Host "after A17:" ($ws.Range("A17").Value2)
